# Auto-generated Excel COM-interop script applying market-price refresh updates
# to the Leve profit tables across multiple job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2872.6924
$ws.Range("I40").Value = 2506.25
$ws.Range("J40").Value = 3459
$ws.Range("K40").Value = 2506.25
$ws.Range("L40").Value = 3459
$ws.Range("M40").Value = -2331.25
$ws.Range("N40").Value = -3809

$ws.Range("H43").Value = 1666
$ws.Range("I43").Value = 1499
$ws.Range("K43").Value = 1499
$ws.Range("M43").Value = -1430

$ws.Range("H137").Value = 2384.5
$ws.Range("I137").Value = 970
$ws.Range("K137").Value = 2910
$ws.Range("M137").Value = -360

$ws.Range("H138").Value = 2628.5696
$ws.Range("I138").Value = 4690.2666
$ws.Range("J138").Value = 2145.3594
$ws.Range("K138").Value = 14070.7998
$ws.Range("L138").Value = 6436.0782
$ws.Range("M138").Value = -8930.799800000001
$ws.Range("N138").Value = -16716.0782

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3248.5483
$ws.Range("I32").Value = 2783.0688
$ws.Range("K32").Value = 2783.0688
$ws.Range("M32").Value = -2496.0688

$ws.Range("H61").Value = 2460.8125
$ws.Range("I61").Value = 2361.2964
$ws.Range("K61").Value = 2361.2964
$ws.Range("M61").Value = -2149.2964

$ws.Range("H74").Value = 1504
$ws.Range("I74").Value = 1504
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1504
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -630
$ws.Range("N74").ClearContents() | Out-Null

$ws.Range("H77").Value = 1504
$ws.Range("I77").Value = 1504
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 7520
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -3152
$ws.Range("N77").Value = -3152

$ws.Range("H136").Value = 2460.8125
$ws.Range("I136").Value = 2361.2964
$ws.Range("K136").Value = 7083.889200000001
$ws.Range("M136").Value = -4533.889200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 9125
$ws.Range("J24").Value = 2500
$ws.Range("L24").Value = 2500
$ws.Range("N24").Value = -2970

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 987.25
$ws.Range("I16").Value = 1153
$ws.Range("K16").Value = 1153
$ws.Range("M16").Value = -866

$ws.Range("H31").Value = 2536.913
$ws.Range("I31").Value = 1755.8889
$ws.Range("K31").Value = 1755.8889
$ws.Range("M31").Value = -1460.8889

$ws.Range("H34").Value = 2536.913
$ws.Range("I34").Value = 1755.8889
$ws.Range("K34").Value = 1755.8889
$ws.Range("M34").Value = -1553.8889

$ws.Range("H58").Value = 2111.739
$ws.Range("J58").Value = 1464.5714
$ws.Range("L58").Value = 1464.5714
$ws.Range("N58").Value = -1870.5714

$ws.Range("H113").Value = 987.25
$ws.Range("I113").Value = 1153
$ws.Range("K113").Value = 1153
$ws.Range("M113").Value = 1017

$ws.Range("H134").Value = 2300.4333
$ws.Range("I134").Value = 2420.25
$ws.Range("J134").Value = 1821.1666
$ws.Range("K134").Value = 7260.75
$ws.Range("L134").Value = 5463.4998
$ws.Range("M134").Value = -4725.75
$ws.Range("N134").Value = -10533.4998

$ws.Range("H136").Value = 2111.739
$ws.Range("J136").Value = 1464.5714
$ws.Range("L136").Value = 4393.7142
$ws.Range("N136").Value = -9493.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 524
$ws.Range("I33").Value = 548.6667
$ws.Range("K33").Value = 3292.0002
$ws.Range("M33").Value = -3009.0002

$ws.Range("H128").Value = 137707.8
$ws.Range("I128").Value = 137707.8
$ws.Range("K128").Value = 413123.4
$ws.Range("M128").Value = -408143.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 500
$ws.Range("J17").Value = 500
$ws.Range("L17").Value = 500
$ws.Range("N17").Value = -836

$ws.Range("H113").Value = 3055.125
$ws.Range("I113").Value = 1924.2
$ws.Range("J113").Value = 4940
$ws.Range("K113").Value = 1924.2
$ws.Range("L113").Value = 4940
$ws.Range("M113").Value = 245.8
$ws.Range("N113").Value = -9280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 9619.916999999999
$ws.Range("J20").Value = 9619.916999999999
$ws.Range("L20").Value = 9619.916999999999
$ws.Range("N20").Value = -10071.917

$ws.Range("H21").Value = 4500
$ws.Range("J21").Value = 4500
$ws.Range("L21").Value = 4500
$ws.Range("N21").Value = -4848

$ws.Range("H43").Value = 3410000
$ws.Range("J43").Value = 6191818
$ws.Range("L43").Value = 6191818
$ws.Range("N43").Value = -6192204

$ws.Range("H46").Value = 3988.5557
$ws.Range("I46").Value = 2724.25
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 2724.25
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -2536.25
$ws.Range("N46").Value = -5376

$ws.Range("H122").Value = 3866.8215
$ws.Range("I122").Value = 4210.294
$ws.Range("J122").Value = 3336
$ws.Range("K122").Value = 12630.882
$ws.Range("L122").Value = 10008
$ws.Range("M122").Value = -10180.882
$ws.Range("N122").Value = -14908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 8000
$ws.Range("J20").Value = 8000
$ws.Range("L20").Value = 8000
$ws.Range("N20").Value = -8480

$ws.Range("H25").Value = 9997
$ws.Range("J25").Value = 9997
$ws.Range("L25").Value = 9997
$ws.Range("N25").Value = -10583

$ws.Range("H32").Value = 42955.25
$ws.Range("J32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("N32").Value = -20634

$ws.Range("H34").Value = 41715
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents() | Out-Null

$ws.Range("H40").Value = 49999
$ws.Range("I40").Value = 49999
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 49999
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -49850
$ws.Range("N40").ClearContents() | Out-Null

$ws.Range("H132").Value = 1412.1154
$ws.Range("I132").Value = 1501.2
$ws.Range("J132").Value = 1115.1666
$ws.Range("K132").Value = 4503.6
$ws.Range("L132").Value = 3345.4998
$ws.Range("M132").Value = -1973.6
$ws.Range("N132").Value = -8405.4998

$ws.Range("H140").Value = 72500
$ws.Range("J140").Value = 72500
$ws.Range("L140").Value = 72500
$ws.Range("N140").Value = -82860
